# Update statistics (mean, median, std, variance) for the btc_stats sheet.
# Columns: B=min, C=max, D=mean, E=median, F=std, G=variance
# Rows:    2=Open, 3=High, 4=Low, 5=Close, 6=Adj Close, 7=Volume

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 15556.57379705767;  E = 8664.668945000001;  F = 17047.82838561126;  G = 290628452.6652529 }
    3 = @{ D = 15982.76269532862;  E = 8836.5161135;        F = 17512.90988451357;  G = 306702012.6230932 }
    4 = @{ D = 15087.92766895919;  E = 8394.783691500001;  F = 16508.60883706253;  G = 272534165.7351391 }
    5 = @{ D = 15576.35411916322;  E = 8664.410156;         F = 17050.26147985278;  G = 290711416.5313514 }
    6 = @{ D = 15576.35411916322;  E = 8664.410156;         F = 17050.26147985278;  G = 290711416.5313514 }
    7 = @{ D = 21091076418.14472; E = 16531574835.5;        F = 21199144177.7595;   G = 449403713869434585088.0 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
